$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "1" to "ამბროლაური"
$ws.Name = "ამბროლაური"

# Remove the census-reference row ("(მოსახლეობის აღწერის შედეგებით)") —
# this was row 2, pushing everything below it up by one.
$ws.Rows(2).Delete()

# The table used to carry three census years (1989 / 2002 / 2014) in
# columns B:D — keep only column B and drop the old 2002/2014 columns
# that used to live in C:D.
$ws.Range("C1:D1048576").EntireColumn.Delete()

# Clear the two stray formatted-but-empty cells in column B that are left
# over from the old header/spacer rows.
$ws.Range("B1").Clear()
$ws.Range("B2").Clear()

# Update the remaining data column so it reflects the latest (2014) figures.
$ws.Range("B4").Value = 2014
$ws.Range("B5").Value = 1139.2
